$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "2aF-3534"
$ws.Range("B3").Value = "asdfasdf"
$ws.Range("C3").Value = -3534
$ws.Range("D3").Value = "Femenino"
$ws.Range("E3").Value = "13/01/2021 - 5:36:10 p. m."
$ws.Range("F3").Value = "asdfasdf"

# Row 4
$ws.Range("A4").Value = "3aF34"
$ws.Range("B4").Value = "asdfqwerwerty"
$ws.Range("C4").Value = 34
$ws.Range("D4").Value = "Femenino"
$ws.Range("E4").Value = "13/01/2021 - 5:44:32 p. m."
$ws.Range("F4").Value = "sdfgwq3456"
